# Fill in the placeholder date on the "Definición del proyecto" cover
# sheet with the actual inspection date.
$d = $word.ActiveDocument

$oldText = "Fecha: dd/mm/aaaa"
$newText = "Fecha: 09/06/2021"

$wdFindStop    = 0   # do not let Find wander past the end of the range
$wdReplaceOne  = 1   # replace only the (single) match inside the range

# The placeholder text appears twice in this checklist template (once on
# the cover block, once again under "Perfiles"); only the first one is
# the project's header date, so scope the Find/Replace to that single
# paragraph's Range instead of the whole story.
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text.Contains($oldText)) {
        $target = $para.Range
        $target.Find.ClearFormatting()
        $target.Find.Execute(
            $oldText,       # FindText
            $true,          # MatchCase
            $false,         # MatchWholeWord
            $false,         # MatchWildcards
            $false,         # MatchSoundsLike
            $false,         # MatchAllWordForms
            $true,          # Forward
            $wdFindStop,    # Wrap
            $false,         # Format
            $newText,       # ReplaceWith
            $wdReplaceOne   # Replace
        )
        break
    }
}
